$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.91"
$ws.Range("E2").Value = "'1.01%"
$ws.Range("G2").Value = "'7"
$ws.Range("D3").Value = "'32.82"
$ws.Range("E3").Value = "'4.42%"
$ws.Range("G3").Value = "'7"
$ws.Range("D4").Value = "'4.960"
$ws.Range("E4").Value = "'-2.63%"
$ws.Range("G4").Value = "'7"
$ws.Range("D5").Value = "'0.07778"
$ws.Range("E5").Value = "'-1.28%"
$ws.Range("G5").Value = "'7"
$ws.Range("D6").Value = "'1.990"
$ws.Range("E6").Value = "'-12.50%"
$ws.Range("G6").Value = "'7"
$ws.Range("D7").Value = "'7.860"
$ws.Range("E7").Value = "'0.72%"
$ws.Range("G7").Value = "'7"
$ws.Range("D8").Value = "'3.800"
$ws.Range("E8").Value = "'-1.43%"
$ws.Range("G8").Value = "'7"
$ws.Range("D9").Value = "'0.9263"
$ws.Range("E9").Value = "'0.70%"
$ws.Range("G9").Value = "'7"
$ws.Range("D10").Value = "'0.1758"
$ws.Range("E10").Value = "'0.76%"
$ws.Range("G10").Value = "'7"
$ws.Range("D11").Value = "'0.07885"
$ws.Range("E11").Value = "'4.02%"
$ws.Range("G11").Value = "'7"
$ws.Range("D12").Value = "'0.08610"
$ws.Range("E12").Value = "'-9.25%"
$ws.Range("G12").Value = "'7"
$ws.Range("D13").Value = "'0.03155"
$ws.Range("E13").Value = "'4.85%"
$ws.Range("G13").Value = "'7"
$ws.Range("D14").Value = "'0.1001"
$ws.Range("E14").Value = "'0.01%"
$ws.Range("G14").Value = "'7"
$ws.Range("D15").Value = "'0.001513"
$ws.Range("E15").Value = "'0.46%"
$ws.Range("G15").Value = "'7"
$ws.Range("D16").Value = "'0.005923"
$ws.Range("E16").Value = "'-2.30%"
$ws.Range("G16").Value = "'7"
$ws.Range("D17").Value = "'3.467"
$ws.Range("E17").Value = "'-0.26%"
$ws.Range("G17").Value = "'7"
$ws.Range("E18").Value = "'-4.02%"
$ws.Range("G18").Value = "'7"
$ws.Range("D19").Value = "'0.3331"
$ws.Range("E19").Value = "'1.81%"
$ws.Range("G19").Value = "'7"
$ws.Range("E20").Value = "'0.73%"
$ws.Range("G20").Value = "'7"
$ws.Range("D21").Value = "'4.341"
$ws.Range("E21").Value = "'10.07%"
$ws.Range("G21").Value = "'7"
$ws.Range("D22").Value = "'0.1993"
$ws.Range("E22").Value = "'16.42%"
$ws.Range("G22").Value = "'7"
$ws.Range("D23").Value = "'0.04562"
$ws.Range("E23").Value = "'-1.25%"
$ws.Range("G23").Value = "'7"
$ws.Range("D24").Value = "'0.001226"
$ws.Range("E24").Value = "'-2.31%"
$ws.Range("G24").Value = "'7"
$ws.Range("E25").Value = "'-0.74%"
$ws.Range("G25").Value = "'7"
$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'0.10%"
$ws.Range("G26").Value = "'7"
$ws.Range("G27").Value = "'7"
$ws.Range("G28").Value = "'7"
$ws.Range("G29").Value = "'7"
$ws.Range("G30").Value = "'7"
$ws.Range("G31").Value = "'7"
$ws.Range("G32").Value = "'7"
$ws.Range("G33").Value = "'7"
$ws.Range("G34").Value = "'7"
$ws.Range("G35").Value = "'7"
$ws.Range("G36").Value = "'7"
$ws.Range("G37").Value = "'7"
$ws.Range("G38").Value = "'7"
$ws.Range("D39").Value = "'0.01709"
$ws.Range("E39").Value = "'-1.30%"
$ws.Range("G39").Value = "'7"
$ws.Range("D40").Value = "'0.04717"
$ws.Range("E40").Value = "'2.43%"
$ws.Range("G40").Value = "'7"
$ws.Range("D41").Value = "'0.007837"
$ws.Range("E41").Value = "'12.99%"
$ws.Range("G41").Value = "'7"
$ws.Range("D42").Value = "'0.1354"
$ws.Range("E42").Value = "'-0.29%"
$ws.Range("G42").Value = "'7"
$ws.Range("D43").Value = "'0.002342"
$ws.Range("E43").Value = "'6.95%"
$ws.Range("G43").Value = "'7"
$ws.Range("D44").Value = "'0.01051"
$ws.Range("E44").Value = "'2.52%"
$ws.Range("G44").Value = "'7"
$ws.Range("D45").Value = "'0.00006277"
$ws.Range("E45").Value = "'0.16%"
$ws.Range("G45").Value = "'7"
$ws.Range("E46").Value = "'0.13%"
$ws.Range("G46").Value = "'7"
$ws.Range("D47").Value = "'0.003104"
$ws.Range("G47").Value = "'7"
$ws.Range("D48").Value = "'0.8234"
$ws.Range("E48").Value = "'10.32%"
$ws.Range("G48").Value = "'7"
$ws.Range("E49").Value = "'0.13%"
$ws.Range("G49").Value = "'7"
$ws.Range("E50").Value = "'0.13%"
$ws.Range("G50").Value = "'7"
$ws.Range("G51").Value = "'7"
